$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the "Mirrorless Camera" product block (rows 57-60) with "Digital Piano".
# Write the query text (column D) before the product name (column B) so that new
# shared-string entries land in the same order as the target workbook (queries
# first, product name last).
$ws.Cells.Item(57, 4).Value = "What is the best digital piano?"
$ws.Cells.Item(58, 4).Value = "What is the best digital piano under 5000 dollars?"
$ws.Cells.Item(59, 4).Value = "What is the best digital piano under 5000 dollars with weighted keys?"
$ws.Cells.Item(60, 4).Value = "What is the best digital piano under 5000 dollars with weighted keys and Bluetooth connectivity?"

$ws.Cells.Item(57, 2).Value = "Digital Piano"
$ws.Cells.Item(58, 2).Value = "Digital Piano"
$ws.Cells.Item(59, 2).Value = "Digital Piano"
$ws.Cells.Item(60, 2).Value = "Digital Piano"

# Replace the "Smart TV" product block (rows 62-65) with "Refrigerator", same
# write order (queries, then product name).
$ws.Cells.Item(62, 4).Value = "What is the best refrigerator?"
$ws.Cells.Item(63, 4).Value = "What is the best refrigerator under 5000 dollars?"
$ws.Cells.Item(64, 4).Value = "What is the best refrigerator under 5000 dollars with a built-in ice maker?"
$ws.Cells.Item(65, 4).Value = "What is the best refrigerator under 5000 dollars with a built-in ice maker and smart Wi-Fi connectivity?"

$ws.Cells.Item(62, 2).Value = "Refrigerator"
$ws.Cells.Item(63, 2).Value = "Refrigerator"
$ws.Cells.Item(64, 2).Value = "Refrigerator"
$ws.Cells.Item(65, 2).Value = "Refrigerator"

# Update the view state: scroll so row 22 is the top visible row, then select H58
# (mirrors the sheetView's topLeftCell="A22" / selection activeCell="H58" in the
# saved workbook).
$ws.Range("A22").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 22
$aw.ScrollColumn = 1
$ws.Range("H58").Select()
